$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# --- Register "Description" column (J) for the register header rows ---
$ws.Range("J4").Value = "ID register"
$ws.Range("J6").Value = "Clock register"
$ws.Range("J11").Value = "Timer register"

# --- Field "Description" column (Q) for the individual field rows ---
$ws.Range("Q5").Value = "ID field"
$ws.Range("Q7").Value = "Clock division"
$ws.Range("Q8").Value = "Clock frequency"
$ws.Range("Q9").Value = "Clock enable"
$ws.Range("Q12").Value = "Counter value"
$ws.Range("Q13").Value = "Timer enable"
$ws.Range("Q14").Value = "Timer start"

# --- EN field's reset value becomes a textual hex literal instead of the number 1 ---
$ws.Range("N9").Value = "0x1"

# --- Move the active selection down one row, matching the resulting file ---
$ws.Range("N14").Select()
